# "1. joystick keymap header 추가"
#
# This edit:
#  1. Fixes the joystick axis label "Ax6" -> "Ax4" in the keymap
#     diagram textbox on slide 1 (the stray label referred to axis 6
#     but should reference axis 4).
#  2. Refreshes the fixed "date updated" footer field (datetimeFigureOut)
#     baked into the slide master and every slide layout from
#     2025-02-03 to 2025-02-10.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Slide 1: "Ax6 [-1 ~ 1] 기본값 : -1" -> "Ax4 [-1 ~ 1] 기본값 : -1"
# ---------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$shapeCount = $slide1.Shapes.Count
for ($i = 1; $i -le $shapeCount; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Length -ge 3) {
            if ($tr.Characters(1, 3).Text -eq "Ax6") {
                # Only touch the exact run that spells out the axis label
                # (leaves the rest of the text box's runs/formatting intact).
                $tr.Characters(1, 3).Text = "Ax4"
            }
        }
    }
}

# ---------------------------------------------------------------
# 2) Master + every layout: fixed footer date 2025-02-03 -> 2025-02-10
# ---------------------------------------------------------------
function Update-FooterDate {
    param($shapes)

    $n = $shapes.Count
    for ($i = 1; $i -le $n; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "2025-02-03") {
                $tr.Text = "2025-02-10"
            }
        }
    }
}

$master = $p.SlideMaster
Update-FooterDate $master.Shapes

$layouts = $master.CustomLayouts
$layoutCount = $layouts.Count
for ($li = 1; $li -le $layoutCount; $li++) {
    Update-FooterDate $layouts.Item($li).Shapes
}
